$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 - this pushes the existing row 4
# (date 2023-02-16 / 44973) down to row 5, preserving all of its values.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44981
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101006
$ws.Range("J4").Value = "Higo"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 3125
$ws.Range("T4").Value = 8

# Match the date-cell number format used by the other rows in column D.
$ws.Range("D4").NumberFormat = $ws.Range("D3").NumberFormat
